$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.994.58'
$ws.Range("E2").Value = '  -3.50%  '
$ws.Range("D3").Value = '3.361.58'
$ws.Range("E3").Value = '  -2.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("E9").Value = '  +1.81%  '
$ws.Range("E10").Value = '  -1.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.417'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.98%  '
$ws.Range("D12").Value = '3.936.39'
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.130'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.08%  '
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("D16").Value = '3.351.61'
$ws.Range("E16").Value = '  -3.02%  '
$ws.Range("D17").Value = '61.057.07'
$ws.Range("E17").Value = '  -3.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("E19").Value = '  -1.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.564'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '3.522.32'
$ws.Range("E25").Value = '  -1.80%  '
$ws.Range("E26").Value = '  -5.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.177'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.31%  '
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.96%  '
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  -5.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.29'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.41'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '169.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("E37").Value = '  -4.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '29.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.55%  '
$ws.Range("D40").Value = '3.396.29'
$ws.Range("E40").Value = '  -2.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0755'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.764'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.66%  '
$ws.Range("E44").Value = '  -1.74%  '
$ws.Range("E45").Value = '  -4.32%  '
$ws.Range("E46").Value = '  -6.47%  '
$ws.Range("D47").Value = '2.517.03'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.84'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("E49").Value = '  -2.15%  '
$ws.Range("E51").Value = '  -2.54%  '
